$wb = $excel.ActiveWorkbook

# --- 1. Update the selection/active state on "Лист1" (sheet 1) -----------
# Target: sheetView no longer tabSelected, selection becomes the whole
# A1:XFD2 block (no single active cell highlighted inside it).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:XFD2").Select() | Out-Null

# --- 2. Add a new worksheet "Лист3" at the end of the workbook -----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Лист3"

# --- 3. Populate Лист3 with the first two rows of Лист1 (same machine) ---
$ws3.Range("A1").Value = "Bernina"
$ws3.Range("B1").Value = "B 570 QE"
$ws3.Range("C1").Value = 5990
$ws3.Range("D1").Value = "BYN"
$ws3.Range("E1").Value = "https://catalog.onliner.by/sewingmachines/bernina/b570qe"
$ws3.Hyperlinks.Add($ws3.Range("E1"), "https://catalog.onliner.by/sewingmachines/bernina/b570qe") | Out-Null
$ws3.Range("E1").Style = "Гиперссылка"

$ws3.Range("A2").Value = "Bernina"
$ws3.Range("B2").Value = "B 570 QE (с вышивальным блоком)"
$ws3.Range("C2").Value = 7700
$ws3.Range("D2").Value = "BYN"
$ws3.Range("E2").Value = "https://catalog.onliner.by/sewingmachines/bernina/b570qe2"
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://catalog.onliner.by/sewingmachines/bernina/b570qe2") | Out-Null
$ws3.Range("E2").Style = "Гиперссылка"

# --- 4. Leave Лист3 as the active sheet/tab with G8 selected --------------
$ws3.Range("G8").Select() | Out-Null
